# TODO.xlsx — mark two more rows as "Done" and scroll the sheet down to
# show the bottom of the list (end-of-game / Chapter 2 teaser rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3  ("Menu Scene")      -> mark Status (col D) as Done
# Row 40 ("Credits scene")   -> mark Status (col D) as Done
$ws.Range("D3").Value = "Done"
$ws.Range("D40").Value = "Done"

# Scroll the view down so row 22 is at the top and select D40, matching
# where the author left the cursor after finishing the list.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D40").Select()
